$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 2.997429241610044

# Row 3 updates
$ws.Range("B3").Value = 0.04763786555579896
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 646.493194727951
